$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C: x - xsr (deviation of each x value from the mean of column A)
$ws.Range("C2").Formula = '=A2-$A$12'
$ws.Range("C3:C11").Formula = '=A3-$A$12'
$ws.Range("C12").Formula = '=AVERAGE(C2:C11)'

$ws.Range("C12").Select()
